$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I4").Value = -0.6627369299571171
$ws.Range("J4").Value = 0.457686881153907
$ws.Range("K4").Value = 0.3379598144967776
$ws.Range("L4").Value = 2.720457785699356
